$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- 1. Fix TestCase_B57 row (row 58): E58 PASS -> SKIP, since a new final
#        "PASS" row is appended at the bottom of the table. ---
$ws.Range("E58").Value = "SKIP"

# --- 2. Append six new "Trending" test case rows (59-64) ---
$ws.Range("A59").Value = "TestCase_B58"
$ws.Range("B59").Value = "OPQA-287"
$ws.Range("C59").Value = "Verify that the following sections get displayed in the search type ahead when user searches using ALL option in the search drop down and that the searched keyword is present in all the sections:`na)Autocompleted keyword`nb)4 suggestions in CATEGORY section`nc)4 suggestions in ARTICLES section`nd)4 suggestions in PATENTS section`ne)4 suggestions in POSTS section`nf)4 suggestions in PEOPLE section"
$ws.Range("D59").Value = "Y"
$ws.Range("E59").Value = "SKIP"

$ws.Range("A60").Value = "TestCase_B59"
$ws.Range("B60").Value = "OPQA-311"
$ws.Range("C60").Value = "Verify that 10 article suggestions get displayed in the search type ahead when user searches using ARTICLES option in the search drop down and that the searched keyword is present in all the suggestions"
$ws.Range("D60").Value = "Y"
$ws.Range("E60").Value = "SKIP"

$ws.Range("A61").Value = "TestCase_B60"
$ws.Range("B61").Value = "OPQA-362"
$ws.Range("C61").Value = "Verify that 10 patent suggestions get displayed in the search type ahead when user searches using PATENTS option in the search drop down and that the searched keyword is present in all the suggestions"
$ws.Range("D61").Value = "Y"
$ws.Range("E61").Value = "SKIP"

$ws.Range("A62").Value = "TestCase_B61"
$ws.Range("B62").Value = "OPQA-371"
$ws.Range("C62").Value = "Verify that 10 post suggestions get displayed in the search type ahead when user searches using POSTS option in the search drop down and that the searched keyword is present in all the suggestions"
$ws.Range("D62").Value = "Y"
$ws.Range("E62").Value = "SKIP"

$ws.Range("A63").Value = "TestCase_B62"
$ws.Range("B63").Value = "OPQA-378"
$ws.Range("C63").Value = "Verify that 10 people suggestions get displayed in the search type ahead when user searches using PEOPLE option in the search drop down and that the searched keyword is present in all the suggestions"
$ws.Range("D63").Value = "Y"
$ws.Range("E63").Value = "SKIP"

$ws.Range("A64").Value = "TestCase_B63"
$ws.Range("B64").Value = "OPQA-258"
$ws.Range("C64").Value = "Verify that no search results get displayed if search engine doesn't interpret the query and that a proper message gets displayed regarding that"
$ws.Range("D64").Value = "Y"
$ws.Range("E64").Value = "PASS"

# --- 3. Apply the same formatting as the row above (TestCase_B57, row 58)
#        to the new rows so borders/fills/wrap match the rest of the table ---
$ws.Range("A58:E58").Copy() | Out-Null
$ws.Range("A59:E64").PasteSpecial(-4122) | Out-Null

# --- 4. Row heights: row 59 holds a 6-line description, rows 60-63 hold a
#        2-line description, row 64 is single-line (default height) ---
$ws.Rows.Item(59).RowHeight = 120
$ws.Rows.Item(60).RowHeight = 30
$ws.Rows.Item(61).RowHeight = 30
$ws.Rows.Item(62).RowHeight = 30
$ws.Rows.Item(63).RowHeight = 30

# --- 5. Normalize column D's style across the whole table (rows 2-64): the
#        workbook re-saved D with the borderId/fillId-only style (like
#        column E) instead of the older "applyFill" variant ---
$ws.Range("E2").Copy() | Out-Null
$ws.Range("D2:D64").PasteSpecial(-4122) | Out-Null
$ws.Range("D2:D64").Value = "Y"

$excel.CutCopyMode = 0

# --- 6. Restore the selected cell to D6 ---
$ws.Range("D6").Select() | Out-Null
